$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell even when the text looks like a
# plain number (e.g. "142.33"), without altering the cell's style (no-op NumberFormat
# round trip) so only the cell value itself changes, matching the source data (which
# stores all of these as text).
function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = '@'
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range('D2').Value = '57.845.26'
$ws.Range('E2').Value = '  +2.34%  '
$ws.Range('D3').Value = '3.069.91'
$ws.Range('E3').Value = '  +2.16%  '
Set-TextValue $ws.Range('D5') '516.87'
$ws.Range('E5').Value = '  +1.75%  '
Set-TextValue $ws.Range('D6') '142.33'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('E11').Value = '  +3.25%  '
$ws.Range('D12').Value = '3.594.84'
$ws.Range('E12').Value = '  +2.35%  '
$ws.Range('E13').Value = '  +3.19%  '
Set-TextValue $ws.Range('D14') '26.17'
$ws.Range('E14').Value = '  +3.26%  '
Set-TextValue $ws.Range('D15') '0.0000164'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = '57.871.33'
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('D17').Value = '3.067.88'
$ws.Range('E17').Value = '  +2.24%  '
Set-TextValue $ws.Range('D18') '6.07'
$ws.Range('E18').Value = '  +2.00%  '
Set-TextValue $ws.Range('D19') '12.87'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('E20').Value = '  +1.36%  '
Set-TextValue $ws.Range('D21') '331.40'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('E22').Value = '  +0.07%  '
Set-TextValue $ws.Range('D23') '0.500'
$ws.Range('E23').Value = '  +0.58%  '
Set-TextValue $ws.Range('D24') '65.72'
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('E25').Value = '  +2.64%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -3.83%  '
Set-TextValue $ws.Range('D28') '6.42'
$ws.Range('E28').Value = '  +1.13%  '
Set-TextValue $ws.Range('D29') '7.21'
$ws.Range('E29').Value = '  +4.76%  '
$ws.Range('E30').Value = '  +1.91%  '
Set-TextValue $ws.Range('D31') '1.19'
$ws.Range('E31').Value = '  +2.73%  '
Set-TextValue $ws.Range('D32') '20.73'
$ws.Range('E32').Value = '  +1.94%  '
Set-TextValue $ws.Range('D33') '154.62'
$ws.Range('E33').Value = '  +0.41%  '
Set-TextValue $ws.Range('D34') '4.53'
$ws.Range('E34').Value = '  +1.71%  '
Set-TextValue $ws.Range('D35') '27.18'
$ws.Range('E35').Value = '  +3.53%  '
Set-TextValue $ws.Range('D36') '5.97'
$ws.Range('E36').Value = '  +2.42%  '
Set-TextValue $ws.Range('D37') '1.27'
$ws.Range('E37').Value = '  +3.33%  '
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = '3.110.69'
$ws.Range('E39').Value = '  +2.42%  '
Set-TextValue $ws.Range('D40') '3.92'
$ws.Range('E40').Value = '  +3.23%  '
Set-TextValue $ws.Range('D41') '36.58'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  +0.07%  '
Set-TextValue $ws.Range('D43') '0.657'
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('D44').Value = '2.256.21'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('E45').Value = '  +8.38%  '
Set-TextValue $ws.Range('D46') '20.74'
$ws.Range('E46').Value = '  +6.34%  '
Set-TextValue $ws.Range('D47') '1.37'
$ws.Range('E47').Value = '  +1.51%  '
Set-TextValue $ws.Range('D48') '0.942'
$ws.Range('E48').Value = '  +1.38%  '
Set-TextValue $ws.Range('D49') '5.90'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('E50').Value = '  +8.49%  '
Set-TextValue $ws.Range('D51') '260.34'
$ws.Range('E51').Value = '  +13.55%  '
